# Updates generated data for gh-pages output (bilibili "想去人数" / attendance
# counters bump, plus a newly scraped event row for 合肥·一生必听的古典系列
# 《钟》—超技钢琴曲炫彩音乐会 on 2024-12-07).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet 1): bump the "想去人数" (F) counters for existing rows.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Cells.Item(2, 6).Value = 655
$wsExpo.Cells.Item(3, 6).Value = 502
$wsExpo.Cells.Item(8, 6).Value = 2065
$wsExpo.Cells.Item(9, 6).Value = 4110
$wsExpo.Cells.Item(10, 6).Value = 97

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet 2): bump F2, then append the new row 4.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)
$wsShow.Cells.Item(2, 6).Value = 59

# Clone formatting of the last existing data row (row 3) onto the new row 4
# so the new row's column-A index cell picks up the same bold/border/
# center-top style used throughout the sheet.
$wsShow.Range("A3:I3").Copy()
$wsShow.Range("A4:I4").PasteSpecial(-4122)

$wsShow.Cells.Item(4, 1).Value = 3

# Column B holds a literal "YYYY-MM-DD" string, not a real date - force text
# formatting before assigning so Excel doesn't auto-convert it to a date
# serial, then reset the cell format (copy plain formatting from D3) so no
# stray number-format style sticks to the cell.
$wsShow.Cells.Item(4, 2).NumberFormat = "@"
$wsShow.Cells.Item(4, 2).Value = "2024-12-07"
$wsShow.Cells.Item(3, 4).Copy()
$wsShow.Cells.Item(4, 2).PasteSpecial(-4122)

$wsShow.Cells.Item(4, 3).Value = "合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会"
$wsShow.Cells.Item(4, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"

# Column E holds a literal time-range string; same text-forcing trick as B.
$wsShow.Cells.Item(4, 5).NumberFormat = "@"
$wsShow.Cells.Item(4, 5).Value = "2024.12.07 19:30-12.07 21:00"
$wsShow.Cells.Item(3, 4).Copy()
$wsShow.Cells.Item(4, 5).PasteSpecial(-4122)

$wsShow.Cells.Item(4, 6).Value = 0
$wsShow.Cells.Item(4, 7).Value = 56
$wsShow.Cells.Item(4, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91608"
$wsShow.Cells.Item(4, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet 4): same bumps as sheet 1 (rows 2,3,8,9,10) plus
# sheet 2's bump (row 11), then append the new row 13 (same new event).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Cells.Item(2, 6).Value = 655
$wsAll.Cells.Item(3, 6).Value = 502
$wsAll.Cells.Item(8, 6).Value = 2065
$wsAll.Cells.Item(9, 6).Value = 4110
$wsAll.Cells.Item(10, 6).Value = 97
$wsAll.Cells.Item(11, 6).Value = 59

$wsAll.Range("A12:I12").Copy()
$wsAll.Range("A13:I13").PasteSpecial(-4122)

$wsAll.Cells.Item(13, 1).Value = 12

$wsAll.Cells.Item(13, 2).NumberFormat = "@"
$wsAll.Cells.Item(13, 2).Value = "2024-12-07"
$wsAll.Cells.Item(12, 4).Copy()
$wsAll.Cells.Item(13, 2).PasteSpecial(-4122)

$wsAll.Cells.Item(13, 3).Value = "合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会"
$wsAll.Cells.Item(13, 4).Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"

$wsAll.Cells.Item(13, 5).NumberFormat = "@"
$wsAll.Cells.Item(13, 5).Value = "2024.12.07 19:30-12.07 21:00"
$wsAll.Cells.Item(12, 4).Copy()
$wsAll.Cells.Item(13, 5).PasteSpecial(-4122)

$wsAll.Cells.Item(13, 6).Value = 0
$wsAll.Cells.Item(13, 7).Value = 56
$wsAll.Cells.Item(13, 8).Value = "https://show.bilibili.com/platform/detail.html?id=91608"
$wsAll.Cells.Item(13, 9).Value = "//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg"
